$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.198361873626709
$ws.Range("B1").Value = 2.11119270324707
$ws.Range("C1").Value = 5.788083553314209
$ws.Range("D1").Value = 0.9916841983795166
$ws.Range("E1").Value = 1.164937019348145
